$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New archetype "Druid" added in column H, row 4
$ws.Range("H4").Value = "Druid"

# "Portal Keeper" removed from U4
$ws.Range("U4").ClearContents()

# Renamed "Dropout Mage" -> "Spellthief"
$ws.Range("J5").Value = "Spellthief"

# Swap E12 and K12 values (Empath <-> ALL)
$ws.Range("E12").Value = "ALL"
$ws.Range("K12").Value = "Empath"

# Update selection to K13
$ws.Range("K13").Select()
